$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21 (pushes the IndianWings rows down to 22-24)
# to add a new EuropeanWings route: Madrid-Barajas (LEMD) -> Berlin-Brandenburg (EDDB)
$ws.Rows(21).Insert()

$ws.Range("A21").Value = "EuropeanWings"
$ws.Range("B21").Value = "Madrid-Barajas"
$ws.Range("D21").Value = "Berlin-Brandenburg"
$ws.Range("E21").Value = "EDDB"
$ws.Range("C21").Value = "LEMD"

# Match the vertical-center alignment style used for newly introduced airport names
$ws.Range("D21").VerticalAlignment = -4108

# Update the active selection to match the post-edit workbook state
$ws.Range("C27").Select() | Out-Null
